$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target range to Text format so numeric-looking strings
# (prices and percentages) are preserved exactly as text, not converted
# to numbers with a numeric style.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "328.59"
$ws.Range("E2").Value = "0.14%"
$ws.Range("D3").Value = "44.20"
$ws.Range("E3").Value = "0.31%"
$ws.Range("D4").Value = "5.404"
$ws.Range("E4").Value = "-1.51%"
$ws.Range("D5").Value = "0.08379"
$ws.Range("E5").Value = "4.18%"
$ws.Range("D6").Value = "1.959"
$ws.Range("E6").Value = "-6.56%"
$ws.Range("D7").Value = "0.9737"
$ws.Range("E7").Value = "2.21%"
$ws.Range("D8").Value = "2.521"
$ws.Range("E8").Value = "-4.24%"
$ws.Range("D9").Value = "0.1135"
$ws.Range("E9").Value = "1.54%"
$ws.Range("D10").Value = "0.1895"
$ws.Range("E10").Value = "0.78%"
$ws.Range("D11").Value = "0.09678"
$ws.Range("E11").Value = "-3.28%"
$ws.Range("D12").Value = "0.04670"
$ws.Range("E12").Value = "-1.73%"
$ws.Range("D13").Value = "0.1062"
$ws.Range("E13").Value = "0.30%"
$ws.Range("D14").Value = "0.001295"
$ws.Range("E14").Value = "1.80%"
$ws.Range("D15").Value = "0.006144"
$ws.Range("E15").Value = "1.99%"
$ws.Range("D16").Value = "3.392"
$ws.Range("E16").Value = "0.59%"
$ws.Range("D17").Value = "4.433"
$ws.Range("E17").Value = "0.40%"
$ws.Range("D18").Value = "0.3329"
$ws.Range("E18").Value = "1.38%"
$ws.Range("D19").Value = "9.145"
$ws.Range("E19").Value = "-9.99%"
$ws.Range("D20").Value = "0.1372"
$ws.Range("E20").Value = "-2.05%"
$ws.Range("D21").Value = "0.2548"
$ws.Range("E21").Value = "2.18%"
$ws.Range("D22").Value = "0.04155"
$ws.Range("E22").Value = "1.54%"
$ws.Range("D23").Value = "0.001296"
$ws.Range("E23").Value = "-1.09%"
$ws.Range("D24").Value = "0.004399"
$ws.Range("E24").Value = "1.25%"
$ws.Range("D25").Value = "0.0001301"
$ws.Range("E25").Value = "3.83%"
$ws.Range("E26").Value = "-20.24%"
$ws.Range("D38").Value = "0.02659"
$ws.Range("E38").Value = "1.32%"
$ws.Range("D39").Value = "0.05643"
$ws.Range("E39").Value = "0.38%"
$ws.Range("D40").Value = "0.007834"
$ws.Range("E40").Value = "2.87%"
$ws.Range("D41").Value = "0.1413"
$ws.Range("E41").Value = "0.88%"
$ws.Range("D42").Value = "0.007365"
$ws.Range("E42").Value = "-0.18%"
$ws.Range("D43").Value = "0.002104"
$ws.Range("E43").Value = "5.81%"
$ws.Range("D44").Value = "0.008634"
$ws.Range("E44").Value = "-2.65%"
$ws.Range("D45").Value = "0.3512"
$ws.Range("D46").Value = "0.00006838"
$ws.Range("E46").Value = "-3.53%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "-0.06%"
$ws.Range("D48").Value = "0.003510"
$ws.Range("E48").Value = "0.26%"
$ws.Range("D49").Value = "0.003533"
$ws.Range("E49").Value = "0.94%"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "-0.06%"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "-0.06%"

# Restore original (unstyled/General) formatting now that values are
# committed as text, so no stray number-format style lingers on cells.
$rng.ClearFormats()
